$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 356 (shifts existing rows 356-378 down to 357-379)
$ws.Rows.Item(356).Insert()

# Populate the newly inserted row 356 with the new weekly price entry
$ws.Range("A356").Value = 4
$ws.Range("B356").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C356").Value = "Los Lagos"
$ws.Range("D356").Value = 44931
$ws.Range("E356").Value = 10
$ws.Range("F356").Value = 100112037
$ws.Range("G356").Value = "Cebollín"
$ws.Range("H356").Value = "Sin especificar"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 70
$ws.Range("K356").Value = 7000
$ws.Range("L356").Value = 7000
$ws.Range("M356").Value = 7000
$ws.Range("N356").Value = "`$/paquete 36 unidades"
$ws.Range("O356").Value = "Región Metropolitana"
$ws.Range("P356").Value = 194
$ws.Range("Q356").Value = 36
$ws.Range("R356").Value = "Hortaliza"
